$wb = $excel.ActiveWorkbook

# "Repayment schedule" sheet: insert a new blank column before the
# existing "Late" column (N), shifting Late/heading/Outstanding one
# column to the right (N->O, O->P, P->Q).
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = 10.7109375

# Make "Repayment schedule" the active sheet/tab and set its selection.
$ws.Activate() | Out-Null
$ws.Range("R6").Select() | Out-Null
